$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 679025.8923849599
$ws.Range("C3").Value = 126789.67668903836
$ws.Range("C4").Value = 64519.75085621432
$ws.Range("C5").Value = 487716.46483970823
$ws.Range("C6").Value = 121243.75723927609
$ws.Range("C7").Value = 171097.71222410485
$ws.Range("C8").Value = 200942.87287758396
$ws.Range("C9").Value = 25238.457391499036
$ws.Range("C10").Value = 160372.4625551049
$ws.Range("C11").Value = 130.6300973915227
